$d = $word.ActiveDocument

# Change 1: fix the typo "requiree deleting." -> "requires deleting."
$d.Content.Find.Execute("requiree deleting.", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "requires deleting.", 2)

# Change 2: "The system display the "Delete another medication?" prompt."
#        -> "The system display the "Exit or delete another medication?" prompt."
$openQuote = [char]0x201C
$closeQuote = [char]0x201D
$oldText = "The system display the " + $openQuote + "Delete another medication?" + $closeQuote + " prompt."
$newText = "The system display the " + $openQuote + "Exit or delete another medication?" + $closeQuote + " prompt."
$d.Content.Find.Execute($oldText, $true, $false, $false, $false, $false, `
                         $true, 1, $false, $newText, 2)
